$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $val) {
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "25.972.71"
$ws.Range("E2").Value = "  -0.17%  "

# Row 3
$ws.Range("D3").Value = "1.643.32"
$ws.Range("E3").Value = "  +0.01%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.004"
$ws.Range("E4").Value = "  -0.26%  "

# Row 5
Set-TextValue $ws.Range("D5") "215.48"
$ws.Range("E5").Value = "  +0.10%  "

# Row 6
$ws.Range("E6").Value = "  -0.28%  "

# Row 7
Set-TextValue $ws.Range("D7") "1.004"
$ws.Range("E7").Value = "  -0.27%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.2576"
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.06425"
$ws.Range("E9").Value = "  +0.35%  "

# Row 10
Set-TextValue $ws.Range("D10") "19.67"
$ws.Range("E10").Value = "  +0.30%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.07775"
$ws.Range("E11").Value = "  +0.55%  "

# Row 12
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D12") "4.272"
$ws.Range("E12").Value = "  +0.45%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.640.98"
$ws.Range("E13").Value = "  -0.16%  "

# Row 14
$ws.Range("D14").Value = "1.869.92"
$ws.Range("E14").Value = "  +0.06%  "

# Row 15
Set-TextValue $ws.Range("D15") "0.5450"
$ws.Range("E15").Value = "  +0.11%  "

# Row 16
$ws.Range("D16").Value = "0.0₅7940"
$ws.Range("E16").Value = "  -0.26%  "

# Row 17
Set-TextValue $ws.Range("D17") "64.51"
$ws.Range("E17").Value = "  +1.27%  "

# Row 18
$ws.Range("D18").Value = "26.009.16"
$ws.Range("E18").Value = "  +0.02%  "

# Row 19
Set-TextValue $ws.Range("D19") "1.004"
$ws.Range("E19").Value = "  -0.34%  "

# Row 20
Set-TextValue $ws.Range("D20") "200.42"
$ws.Range("E20").Value = "  -2.63%  "

# Row 21
Set-TextValue $ws.Range("D21") "4.385"
$ws.Range("E21").Value = "  +0.75%  "

# Row 22
Set-TextValue $ws.Range("D22") "9.920"
$ws.Range("E22").Value = "  -0.72%  "

# Row 23
Set-TextValue $ws.Range("D23") "5.983"
$ws.Range("E23").Value = "  +0.11%  "

# Row 24
Set-TextValue $ws.Range("D24") "1.005"
$ws.Range("E24").Value = "  -0.25%  "

# Row 25
Set-TextValue $ws.Range("D25") "1.885"
$ws.Range("E25").Value = "  -3.32%  "

# Row 26
Set-TextValue $ws.Range("D26") "140.69"
$ws.Range("E26").Value = "  -1.48%  "

# Row 27
Set-TextValue $ws.Range("D27") "0.1138"
$ws.Range("E27").Value = "  -1.86%  "

# Row 28
Set-TextValue $ws.Range("D28") "6.835"
$ws.Range("E28").Value = "  -0.43%  "

# Row 29
$ws.Range("E29").Value = "  -0.43%  "

# Row 30
Set-TextValue $ws.Range("D30") "1.244"
$ws.Range("E30").Value = "  +0.58%  "

# Row 31
Set-TextValue $ws.Range("D31") "0.04922"
$ws.Range("E31").Value = "  -2.12%  "

# Row 32
Set-TextValue $ws.Range("D32") "3.271"
$ws.Range("E32").Value = "  -0.67%  "

# Row 33
Set-TextValue $ws.Range("D33") "3.217"
$ws.Range("E33").Value = "  +0.16%  "

# Row 34
Set-TextValue $ws.Range("D34") "1.544"
$ws.Range("E34").Value = "  +0.36%  "

# Row 35
Set-TextValue $ws.Range("D35") "2.372"
$ws.Range("E35").Value = "  +1.25%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.8945"
$ws.Range("E36").Value = "  -1.68%  "

# Row 37
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "1.156.76"
$ws.Range("E37").Value = "  +2.19%  "

# Row 38
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D38") "2.606"
$ws.Range("E38").Value = "  -1.39%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.5585"
$ws.Range("E39").Value = "  -1.49%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.01570"
$ws.Range("E40").Value = "  +0.48%  "

# Row 41
Set-TextValue $ws.Range("D41") "1.004"
$ws.Range("E41").Value = "  -0.28%  "

# Row 42
Set-TextValue $ws.Range("D42") "5.729"
$ws.Range("E42").Value = "  +2.04%  "

# Row 43
Set-TextValue $ws.Range("D43") "0.8129"
$ws.Range("E43").Value = "  -0.72%  "

# Row 44
Set-TextValue $ws.Range("D44") "99.82"
$ws.Range("E44").Value = "  +0.10%  "

# Row 45
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.0₈120"
$ws.Range("E45").Value = "  +3.29%  "

# Row 46
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.780.49"
$ws.Range("E46").Value = "  -0.02%  "

# Row 47
Set-TextValue $ws.Range("D47") "0.4518"
$ws.Range("E47").Value = "  -0.21%  "

# Row 48
Set-TextValue $ws.Range("D48") "1.006"
$ws.Range("E48").Value = "  +0.09%  "

# Row 49
Set-TextValue $ws.Range("D49") "54.89"
$ws.Range("E49").Value = "  -0.04%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.05051"
$ws.Range("E50").Value = "  +0.00%  "

# Row 51
Set-TextValue $ws.Range("D51") "1.004"
$ws.Range("E51").Value = "  -0.15%  "
